$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("applicator")

$ws.Range("Z1000").Value = "XYZPLACEHOLDER"
$found = $ws.Cells.Replace("XYZPLACEHOLDER", "05/03/2018")
Write-Host "replace result: $found"
$v = $ws.Range("Z1000").Value2
Write-Host "Z1000 value2=$v"

$ws.Range("Z1001").Value = "XYZPLACEHOLDER2"
$found2 = $ws.Cells.Replace("XYZPLACEHOLDER2", "**")
Write-Host "replace2 result: $found2"
$v2 = $ws.Range("Z1001").Value2
Write-Host "Z1001 value2=$v2"
